# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# de-de handback has completed and the zh-cn / de-de target + handback
# files / dates are now populated, and the overall status string changes
# from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b89a2ca6c94658e1b717444881e2369ea07d1e8c/e2e/24f3fbdb-43db-46bb-a4aa-3a64f07f679e.md"
$mdDisplay = "24f3fbdb-43db-46bb-a4aa-3a64f07f679e.md"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update the "Status" text everywhere it is used (it is the same
#    shared string reused on the Overview sheet as well as on each of
#    the per-locale sheets), changing it from "Ready for handoff" to
#    "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in the Latest Target File / Latest Handback File /
#    Latest Handback DateTime columns (I, J, K) for row 2.
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = $mdDisplay
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdDisplay)
$zhcn.Range("I2").Font.Name = "Calibri"
$zhcn.Range("I2").Font.Size = 11
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276

$zhcn.Range("J2").Value = "24f3fbdb-43db-46bb-a4aa-3a64f07f679e.02db8b92cf30802664081aa8dbe6dc337d4cbd24.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-05 21:08:43"

# ---------------------------------------------------------------------
# 3. de-de sheet: fill in the Latest Target File / Latest Handback File /
#    Latest Handback DateTime columns (I, J, K) for row 2.
# ---------------------------------------------------------------------
$dede.Range("I2").Value = $mdDisplay
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdDisplay)
$dede.Range("I2").Font.Name = "Calibri"
$dede.Range("I2").Font.Size = 11
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276

$dede.Range("J2").Value = "24f3fbdb-43db-46bb-a4aa-3a64f07f679e.02db8b92cf30802664081aa8dbe6dc337d4cbd24.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 21:08:51"

# ---------------------------------------------------------------------
# 4. Re-fit the column widths that are affected by the longer text that
#    was just written (Status column, Latest Target File, Latest
#    Handback File). The ColumnWidth setter expects "characters" rather
#    than the raw stored column width, so subtract the fixed 5/6
#    padding that the engine re-adds when it persists the sheet.
# ---------------------------------------------------------------------
$padding = 5.0 / 6.0

$overview.Range("E1").ColumnWidth = 29.9777047293527 - $padding
$overview.Range("F1").ColumnWidth = 29.9777047293527 - $padding

$zhcn.Range("C1").ColumnWidth = 29.9777047293527 - $padding
$zhcn.Range("I1").ColumnWidth = 40 - $padding
$zhcn.Range("J1").ColumnWidth = 40 - $padding

$dede.Range("C1").ColumnWidth = 29.9777047293527 - $padding
$dede.Range("I1").ColumnWidth = 40 - $padding
$dede.Range("J1").ColumnWidth = 40 - $padding

Write-Host "Handback report generated."
